$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the 703-765 data block, shifting the
# existing rows (old 703-765) down to (705-767), matching the diff.
$ws.Range("A703:R704").EntireRow.Insert()

# New row 703: Primera quality, week of 2023-06-29 (serial 45106)
$ws.Range("A703").Value = 3
$ws.Range("B703").Value = "Femacal de La Calera"
$ws.Range("C703").Value = "Coquimbo"
$ws.Range("D703").Value = 45106
$ws.Range("E703").Value = 5
$ws.Range("F703").Value = 100112037
$ws.Range("G703").Value = "Cebollín"
$ws.Range("H703").Value = "Sin especificar"
$ws.Range("I703").Value = "Primera"
$ws.Range("J703").Value = 238
$ws.Range("K703").Value = 3800
$ws.Range("L703").Value = 4000
$ws.Range("M703").Value = 3899
$ws.Range("N703").Value = "$/paquete 36 unidades"
$ws.Range("O703").Value = "Provincia de Quillota"
$ws.Range("P703").Value = 108
$ws.Range("Q703").Value = 36
$ws.Range("R703").Value = "Hortaliza"

# New row 704: Segunda quality, same week (serial 45106)
$ws.Range("A704").Value = 3
$ws.Range("B704").Value = "Femacal de La Calera"
$ws.Range("C704").Value = "Coquimbo"
$ws.Range("D704").Value = 45106
$ws.Range("E704").Value = 5
$ws.Range("F704").Value = 100112037
$ws.Range("G704").Value = "Cebollín"
$ws.Range("H704").Value = "Sin especificar"
$ws.Range("I704").Value = "Segunda"
$ws.Range("J704").Value = 130
$ws.Range("K704").Value = 3000
$ws.Range("L704").Value = 3000
$ws.Range("M704").Value = 3000
$ws.Range("N704").Value = "$/paquete 36 unidades"
$ws.Range("O704").Value = "Provincia de Quillota"
$ws.Range("P704").Value = 83
$ws.Range("Q704").Value = 36
$ws.Range("R704").Value = "Hortaliza"
